$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the discontinued set (row 4: 10357 "Shelby Cobra 427 S/C").
# Deleting the entire row shifts all subsequent rows up by one.
$ws.Rows.Item(4).Delete()

# The "N/A" placeholders in the Collection column (D) are being cleared out.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq "N/A") {
        $cell.Value2 = ""
    }
}
